$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 3143.3845
$ws.Range("I33").Value = 2885.6667
$ws.Range("J33").Value = 3723.25
$ws.Range("K33").Value = 2885.6667
$ws.Range("L33").Value = 3723.25
$ws.Range("M33").Value = -2656.6667
$ws.Range("N33").Value = -4181.25
$ws.Range("H70").Value = 2666
$ws.Range("I70").Value = 1999
$ws.Range("K70").Value = 5997
$ws.Range("M70").Value = -5727
$ws.Range("H73").Value = 2666
$ws.Range("I73").Value = 1999
$ws.Range("K73").Value = 5997
$ws.Range("M73").Value = -5061
$ws.Range("H92").Value = 519.4545000000001
$ws.Range("I92").Value = 471.5
$ws.Range("K92").Value = 471.5
$ws.Range("M92").Value = 776.5
$ws.Range("H99").Value = 568.5454999999999
$ws.Range("I99").Value = 569.375
$ws.Range("K99").Value = 1708.125
$ws.Range("M99").Value = -210.125
$ws.Range("H125").Value = 3966.3333
$ws.Range("J125").Value = 3966.3333
$ws.Range("L125").Value = 35696.9997
$ws.Range("N125").Value = -40616.9997
$ws.Range("H137").Value = 4176.643
$ws.Range("I137").Value = 3039.1667
$ws.Range("K137").Value = 9117.500100000001
$ws.Range("M137").Value = -6567.500100000001
$ws.Range("H138").Value = 2647.6438
$ws.Range("I138").Value = 1699.8
$ws.Range("J138").Value = 2705.439
$ws.Range("K138").Value = 5099.4
$ws.Range("L138").Value = 8116.316999999999
$ws.Range("M138").Value = 40.60000000000036
$ws.Range("N138").Value = -18396.317

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12851095
$ws.Range("I32").Value = 20873838
$ws.Range("K32").Value = 20873838
$ws.Range("M32").Value = -20873551
$ws.Range("H61").Value = 37504052
$ws.Range("I61").Value = 29415354
$ws.Range("K61").Value = 29415354
$ws.Range("M61").Value = -29415142
$ws.Range("H74").Value = 13165338
$ws.Range("I74").Value = 19231800
$ws.Range("K74").Value = 19231800
$ws.Range("M74").Value = -19230926
$ws.Range("H77").Value = 13165338
$ws.Range("I77").Value = 19231800
$ws.Range("K77").Value = 96159000
$ws.Range("M77").Value = -96154632
$ws.Range("H109").Value = 52326.332
$ws.Range("J109").Value = 52326.332
$ws.Range("L109").Value = 52326.332
$ws.Range("N109").Value = -55100.332
$ws.Range("H110").Value = 1450.15
$ws.Range("I110").Value = 1201.7858
$ws.Range("K110").Value = 1201.7858
$ws.Range("M110").Value = 843.2141999999999
$ws.Range("H132").Value = 4822.579
$ws.Range("I132").Value = 3868.2778
$ws.Range("J132").Value = 22000
$ws.Range("K132").Value = 11604.8334
$ws.Range("L132").Value = 66000
$ws.Range("M132").Value = -9074.8334
$ws.Range("N132").Value = -71060
$ws.Range("H136").Value = 37504052
$ws.Range("I136").Value = 29415354
$ws.Range("K136").Value = 88246062
$ws.Range("M136").Value = -88243512

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 40699.89
$ws.Range("J81").Value = 40699.89
$ws.Range("L81").Value = 40699.89
$ws.Range("N81").Value = -42821.89
$ws.Range("H84").Value = 40699.89
$ws.Range("J84").Value = 40699.89
$ws.Range("L84").Value = 122099.67
$ws.Range("N84").Value = -132707.67
$ws.Range("H107").Value = 4129.6665
$ws.Range("I107").Value = 4281.7856
$ws.Range("K107").Value = 4281.7856
$ws.Range("M107").Value = -2361.7856
$ws.Range("H122").Value = 60472.5
$ws.Range("J122").Value = 60472.5
$ws.Range("L122").Value = 60472.5
$ws.Range("N122").Value = -70272.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 872597.5
$ws.Range("I31").Value = 11145.134
$ws.Range("K31").Value = 11145.134
$ws.Range("M31").Value = -10850.134
$ws.Range("H34").Value = 872597.5
$ws.Range("I34").Value = 11145.134
$ws.Range("K34").Value = 11145.134
$ws.Range("M34").Value = -10943.134
$ws.Range("H58").Value = 7549.143
$ws.Range("I58").Value = 6866
$ws.Range("K58").Value = 6866
$ws.Range("M58").Value = -6663
$ws.Range("H86").Value = 89502.5
$ws.Range("I86").Value = 6680.8
$ws.Range("J86").Value = 148660.86
$ws.Range("K86").Value = 6680.8
$ws.Range("L86").Value = 148660.86
$ws.Range("M86").Value = -5557.8
$ws.Range("N86").Value = -150906.86
$ws.Range("H89").Value = 89502.5
$ws.Range("I89").Value = 6680.8
$ws.Range("J89").Value = 148660.86
$ws.Range("K89").Value = 33404
$ws.Range("L89").Value = 743304.2999999999
$ws.Range("M89").Value = -27788
$ws.Range("N89").Value = -754536.2999999999
$ws.Range("H94").Value = 3193.2173
$ws.Range("I94").Value = 2398.8572
$ws.Range("J94").Value = 3540.75
$ws.Range("K94").Value = 2398.8572
$ws.Range("L94").Value = 3540.75
$ws.Range("M94").Value = -1947.8572
$ws.Range("N94").Value = -4442.75
$ws.Range("H107").Value = 1700.5
$ws.Range("I107").Value = 1300.875
$ws.Range("K107").Value = 1300.875
$ws.Range("M107").Value = 619.125
$ws.Range("H122").Value = 1411.1818
$ws.Range("I122").Value = 1403.3
$ws.Range("J122").Value = 1490
$ws.Range("K122").Value = 4209.9
$ws.Range("L122").Value = 4470
$ws.Range("M122").Value = -1759.9
$ws.Range("N122").Value = -9370
$ws.Range("H132").Value = 6072.773
$ws.Range("I132").Value = 2419.3333
$ws.Range("K132").Value = 7257.999899999999
$ws.Range("M132").Value = -4727.999899999999
$ws.Range("H136").Value = 7549.143
$ws.Range("I136").Value = 6866
$ws.Range("K136").Value = 20598
$ws.Range("M136").Value = -18048

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 160.38235
$ws.Range("I2").Value = 78.333336
$ws.Range("K2").Value = 470.000016
$ws.Range("M2").Value = -357.000016
$ws.Range("H38").Value = 351.5
$ws.Range("I38").Value = 484.6
$ws.Range("J38").Value = 129.66667
$ws.Range("K38").Value = 1453.8
$ws.Range("L38").Value = 389.00001
$ws.Range("M38").Value = -1106.8
$ws.Range("N38").Value = -1083.00001
$ws.Range("H107").Value = 984.4167
$ws.Range("J107").Value = 1131.3
$ws.Range("L107").Value = 3393.9
$ws.Range("N107").Value = -7233.9
$ws.Range("H122").Value = 2547.963
$ws.Range("I122").Value = 460
$ws.Range("J122").Value = 3022.5
$ws.Range("K122").Value = 4140
$ws.Range("L122").Value = 27202.5
$ws.Range("M122").Value = -1690
$ws.Range("N122").Value = -32102.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1789.8182
$ws.Range("I97").Value = 1789.8182
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1789.8182
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1293.8182
$ws.Range("N97").ClearContents()
$ws.Range("H122").Value = 2383.2144
$ws.Range("I122").Value = 2351.1538
$ws.Range("K122").Value = 7053.4614
$ws.Range("M122").Value = -4603.4614
$ws.Range("H136").Value = 13298.667
$ws.Range("J136").Value = 13298.667
$ws.Range("L136").Value = 39896.001
$ws.Range("N136").Value = -44996.001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1093.5
$ws.Range("I16").Value = 919.4286
$ws.Range("K16").Value = 919.4286
$ws.Range("M16").Value = -749.4286
$ws.Range("H61").Value = 3298.6667
$ws.Range("J61").Value = 2599.6667
$ws.Range("L61").Value = 2599.6667
$ws.Range("N61").Value = -3003.6667
$ws.Range("H68").Value = 3944.5
$ws.Range("I68").Value = 3944.5
$ws.Range("K68").Value = 3944.5
$ws.Range("M68").Value = -3195.5
$ws.Range("H71").Value = 3944.5
$ws.Range("I71").Value = 3944.5
$ws.Range("K71").Value = 19722.5
$ws.Range("M71").Value = -15978.5
$ws.Range("H100").Value = 4997.5
$ws.Range("I100").Value = 3996.6667
$ws.Range("J100").Value = 8000
$ws.Range("K100").Value = 3996.6667
$ws.Range("L100").Value = 8000
$ws.Range("M100").Value = -3455.6667
$ws.Range("N100").Value = -9082
$ws.Range("H113").Value = 3298.6667
$ws.Range("J113").Value = 2599.6667
$ws.Range("L113").Value = 2599.6667
$ws.Range("N113").Value = -6939.6667
$ws.Range("H123").Value = 60000
$ws.Range("J123").Value = 60000
$ws.Range("L123").Value = 60000
$ws.Range("N123").Value = -69800

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 95996.664
$ws.Range("I95").Value = 0
$ws.Range("K95").Value = 0
$ws.Range("M95").ClearContents()
$ws.Range("H100").Value = 546.2917
$ws.Range("I100").Value = 438.22223
$ws.Range("K100").Value = 876.44446
$ws.Range("M100").Value = -335.44446
$ws.Range("H107").Value = 26317080
$ws.Range("I107").Value = 31251120
$ws.Range("K107").Value = 93753360
$ws.Range("M107").Value = -93751440
$ws.Range("H126").Value = 2099.25
$ws.Range("I126").Value = 2099.25
$ws.Range("K126").Value = 6297.75
$ws.Range("M126").Value = -3827.75
